$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing data range before rewriting it with the new row order/content
$ws.Range("A1:B25").ClearContents()

$ws.Range("A1").Value = 'AH = High Strength Steel, AS = Standard Strength Steel,  AX = Extra High Strength Steel, AT = HSA840 Steel, ATR = AT steel with High Rib, AXR = AX steel with High Rib '
$ws.Range("B1").Value = 'Bolt Matrix'

$ws.Range("A2").Value = 'DS = Double ended bolt Standard Strength, CS = Solid Deformed Bar, DCS = Double Corrosion Protection Steel bell'
$ws.Range("B2").Value = 'Bolt Matrix'

$ws.Range("A3").Value = 'DCS = Double Corrosion Protection Plastic bell,  DGB = VB Dynamic Bolt, DS = Double Ended Expansion Shell Bolt, EB = Eye bolt '
$ws.Range("B3").Value = 'Bolt Matrix'

$ws.Range("A4").Value = 'PCS = Paddle Bolt D/Bar, PH = Paddle Bolt, WS = Wiggle Bolt, WT = Threadbar Wiggle Bolt, HTB = Hollow Tube Bolt'
$ws.Range("B4").Value = 'Bolt Matrix'

$ws.Range("A5").Value = 'G = Galvanised, TD = Thermal Diffused, TS = Thermal Metal Spray, P = Plas coating'
$ws.Range("B5").Value = 'Bolt Matrix'

$ws.Range("A6").Value = 'WA =  Perth, B = Brisbane, no suffix = Newcastle'
$ws.Range("B6").Value = 'Bolt Matrix'

$ws.Range("A7").Value = 'EULA = EziTen 21.8mm plain Ultra Strand, EULN = EziTen 21.8mm Indented Ultra Strand'
$ws.Range("B7").Value = 'Strand Matrix'

$ws.Range("A8").Value = 'HTC8N = 28mm Titan Hollow Indented Cable, HTC8P = 28mm Titan Hollow Plain Cable'
$ws.Range("B8").Value = 'Strand Matrix'

$ws.Range("A9").Value = 'ULA = 21.8mm Ultra Strand, ULN = 21.8mm Indented Ultra Strand Cable'
$ws.Range("B9").Value = 'Strand Matrix'

$ws.Range("A10").Value = 'CB = Cable Bolt, CABOLT = Coil Pack'
$ws.Range("B10").Value = 'Cable Matrix'

$ws.Range("A11").Value = 'G = Galvanised, DB = Debonding Tube, GT = Grout Tube, BT =Breather Tube, RD = Refer Drawing, Q =Steel Aglet'
$ws.Range("B11").Value = 'Cable Matrix'

$ws.Range("A12").Value = 'WA = Freight Cost, B = Brisbane, N = Sea freight'
$ws.Range("B12").Value = 'Cable Matrix'

$ws.Range("A13").Value = 'FG = Glass Reinforced Plastic, P = Plastic                  '
$ws.Range("B13").Value = 'DYWI-GRiP'

$ws.Range("A14").Value = 'B = Bolt without continuous thread, D =  Dowel with continuous thread, DH = Dowel Hollow '
$ws.Range("B14").Value = 'DYWI-GRiP'

$ws.Range("A15").Value = 'DD = Dome plate double, D = Dome plate, FCB = Flat with extra hole, F = Flat, HAD = High deflection angle, STP = Star plate '
$ws.Range("B15").Value = 'Plate Matrix'

$ws.Range("A16").Value = 'TDD = Turtle or Jelly mould or Meshing, STX = Extra form, OCP =Octo plate'
$ws.Range("B16").Value = 'Plate Matrix'

$ws.Range("A17").Value = 'WA = Freight included to Perth, B=Brisbane'
$ws.Range("B17").Value = 'Plate Matrix'

$ws.Range("A18").Value = 'G = Galvanised, TD = Thermal Diffused, TS = Thermal Metal Spray, P = Plas coating, R = Rumbled, XB = Bundled for Galvanisining '
$ws.Range("B18").Value = 'Plate Matrix'

$ws.Range("A19").Value = 'FL = Flared hole, N = Indicator Posts, SS316= Stainless Steel'
$ws.Range("B19").Value = 'Plate Matrix'

$ws.Range("A20").Value = 'FB =Friction Bolt, FBS = Friction Bolt Strengthened Ring, FBSP = Friction Bolt Spiral, OM24 = Omega bolt 24 tonne, FBMC=Kinloc SE, FBMD=Kinloc INDIE'
$ws.Range("B20").Value = 'Friction Bolts'

$ws.Range("A21").Value = 'FBSPD = Friction Bolt Spiral large Dome								'
$ws.Range("B21").Value = 'Friction Bolts'

$ws.Range("A22").Value = 'G=Gal, TD=Thermal Diffised, TSP = Thermal Spray Plas coating'
$ws.Range("B22").Value = 'Friction Bolts'

$ws.Range("A23").Value = 'WB = Weld Mesh Bright, WG = Weldmesh Galvanised, WSS = Weld Mesh Stainless Steel, MMB = Mesh Modules Bright, MMG = Mesh Module galvanised, MMSB = Mesh Module Strap Bright, MMSG = Mesh Module Strap Galvanised, MMBB = Mesh Module Bent Bright, MLB=Mesh L Bend '
$ws.Range("B23").Value = 'Mesh Matrix'

$ws.Range("A24").Value = 'B=Brisbane, WA = Western Australia, EXP=Export'
$ws.Range("B24").Value = 'Mesh Matrix'

$ws.Range("A25").Value = 'S=Spanner'
$ws.Range("B25").Value = 'Dollies & Spanners Coal'

$ws.Range("A26").Value = 'PD22 = Pixi Drive 22mm Hex,   SD22RE = Square Drive 22mm Round Extended section above square for chucks with retaining ring'
$ws.Range("B26").Value = 'Dollies & Spanners Coal'

$ws.Range("A27").Value = 'S=Spanner'
$ws.Range("B27").Value = 'Dollies & Spanners Hardrock'

$ws.Range("A28").Value = 'FB39 = Friction bolt dolly 39mm, FB47 = Friction bolt dolly 47mm, IT = Tapered Drive, PR22 = Pixi Round 22mm drive, R38 = Rope thread 38mm, T38 = Tapezoidal thread 38mm'
$ws.Range("B28").Value = 'Dollies & Spanners Hardrock'

$ws.Range("A29").Value = 'R38F = Rope thread 38mm female socket, T38F = Tapezoidal thread 38mm female socket, FB33 = Friction bolt dolly 33mm, FB39 = Friction bolt dolly 39mm, H36F = Hexagonal 36mm female socket'
$ws.Range("B29").Value = 'Dollies & Spanners Hardrock'

$ws.Range("A30").Value = 'WA = Western Australia, B = Brisbane'
$ws.Range("B30").Value = 'Dollies & Spanners Hardrock'

$ws.Range("A31").Value = 'RA = Resin Anchor, RO = Resin Oil Base, RW = Resin Water Base, RL = Resin Low Viscosity, RT = Resin Twin Speed, RTX= Resin Twin speed Extreme, RS = Resin Single Speed, RSX=Resin Single Extreme'
$ws.Range("B31").Value = 'Resin anchors'

$ws.Range("A32").Value = 'WA = Western Australia, B = Brisbane, E=Export, HP = Half Pallet, R = Exworks Rocbolt, D=Dunnage bags'
$ws.Range("B32").Value = 'Resin anchors'

$ws.Range("A33").Value = 'SB = Spade Bit, EN = Eccentric Tip Negative Rake, KS = Modified Spade, EP = Eccentric Tip Positive rake, RB = Rib Drilling Eccentric Tip Positive Rake, EP = Eccentric tip positive Rake, RE = Two wing coal Bit Eccentric tip positive rake, CN = Concentric tip negative rake, PDC = Polycrystalline Diamond Compact  '
$ws.Range("B33").Value = 'Bits for catalogue'

$ws.Range("A34").Value = 'MH = Mine Hanger'
$ws.Range("B34").Value = 'Mine Hangers'

$ws.Range("A35").Value = '74 = Keyhole Bracket, 75 = Keyhole Bracket, 79 = Nuts,  80 = Hooks'
$ws.Range("B35").Value = 'Mine Hangers'

$ws.Range("A36").Value = '7948 = Closed loop in 12mm round bar to suit 20mm Thread Bar (36AF Nut), 7985 = Closed Loop 12mm Round Bar to suit M24 thread (OZ nut 36AF), 8033 = 3 x Open Hooks in 16mm round bar to suit 110mm pipe, 8036 = 4 x Open Hooks in 16mm round bar to suit 110mm pipe, 7471 = Single keyhole bracket with 26mm mounting hole, 7513 = Fan hanger "T" shaped in 16 & 12mm plate with a single keyhole, 7514 = 6 x Keyholes with 3 x 28mm mounting holes, 7523 = 3 Keyholes and 8mm chains with 3 x 26mm mounting holes, 7525 = 1 Keyhole and 8mm chain with 1 x 26mm mounting hole, 7526 = 7 Keyholes with 3 x 26mm mounting holes'
$ws.Range("B36").Value = 'Mine Hangers'

$ws.Range("A37").Value = 'G = Galvanised, EZ = Electroplated Zinc'

# Column A width changed from 252.6640625 (bestFit) to 255.6640625 (explicit, not bestFit)
$ws.Columns.Item(1).ColumnWidth = 254.8333333

# Update the active selection to match the author's final view state
$ws.Range("A41").Select()

